# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload" update
# Target sheet is MAR-22 (already the active sheet / tab 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: tidy the wording of the "Accounting statements" comment
#     (drop the mid-sentence line-break so "validate" and "data" sit on one
#     line) and shrink the row now that it only wraps to two lines.
$ws.Range("D32").Value = "2. Accounting statements has been generated for the GL of the Jan22 for all three centers and shared to Rahman san to validate data"
$ws.Rows.Item(32).RowHeight = 28.8

# --- Row 33: new daily entry (#15, 21-Mar-2022, RPA GSS)
$ws.Range("A31").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = 15

$ws.Range("B31").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B33").Value = 44641

$ws.Range("C31").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = "RPA GSS"

$ws.Range("D31").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = "1. Correction Received for the Warranty Daily task due to no data found issue at SAMSUNG site for the all SSCs, `nit has been fixed, tested in all SSCs and it is running smoothly"

$ws.Range("E31").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = 1

$ws.Range("F31").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value = "Completed"

$ws.Rows.Item(33).RowHeight = 28.8

# --- Row 34: continuation bullet #2 for the same entry
$ws.Range("D34").Value = "2. Public holidays has been implemented at Warranty, Activity and GRS_Details tasks, it has been tested and it is running smoothly."

$ws.Range("E31").Copy()
$ws.Range("E34").PasteSpecial(-4122)
$ws.Range("E34").Value = 1

$ws.Range("F34").Value = "Completed"

# --- Row 35: continuation bullet #3 for the same entry
$ws.Range("D35").Value = "3. In Rlogic, no data failure issue  at warranty task due to samsung site, and it has been fixed, tested and it is running smoothly"

$ws.Range("E31").Copy()
$ws.Range("E35").PasteSpecial(-4122)
$ws.Range("E35").Value = 1

$ws.Range("F35").Value = "Completed"

# --- Update the saved selection / scroll position to match where the
#     author ended up after entering the new rows.
$ws.Range("D35").Select()
